$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Paragraph formatting: add a paragraph border (5pt spacing on each side)
# and widen the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5
$p1.LeftIndent = 11.25

# Remove the trailing run that only contains a single space character.
$trailingSpace = $p1.Range.Duplicate
$trailingSpace.Start = $trailingSpace.End - 2
$trailingSpace.End = $trailingSpace.End - 1
$trailingSpace.Delete()

# Update the placeholder id text in the remaining run.
$null = $d.Content.Find.Execute("**ID__AFFARS_5337_topic_11__ID**", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "**ID__AFFARS_SUBPART_5337_74__ID**", 2)
